$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain text / percentages / multi-dot numbers
# These can be assigned directly without risk of numeric auto-conversion.
$ws.Range("D2").Value = "25.794.12"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "1.637.22"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("E10").Value = "  +4.55%  "
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.638.74"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.865.12"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("D16").Value = "0.0₅7636"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").Value = "25.826.64"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("E20").Value = "  -0.95%  "
$ws.Range("E21").Value = "  -0.91%  "
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("E23").Value = "  +1.35%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  -4.62%  "
$ws.Range("E26").Value = "  -0.73%  "
$ws.Range("E27").Value = "  -1.64%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("E32").Value = "  +1.34%  "
$ws.Range("E33").Value = "  +1.30%  "
$ws.Range("E34").Value = "  +1.76%  "
$ws.Range("E35").Value = "  +0.60%  "
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("D39").Value = "1.130.75"
$ws.Range("E39").Value = "  +1.35%  "
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("E43").Value = "  +0.84%  "
$ws.Range("E44").Value = "  +1.41%  "
$ws.Range("D45").Value = "1.777.00"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("D46").Value = "0.0₈110"
$ws.Range("E46").Value = "  -5.78%  "
$ws.Range("E47").Value = "  +1.58%  "
$ws.Range("E48").Value = "  -3.91%  "
$ws.Range("E49").Value = "  +3.06%  "
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("E51").Value = "  -0.30%  "

# Cells whose new values look like plain numbers (e.g. "0.5056"). Excel would
# otherwise auto-convert these to numeric cells, but the source data stores them
# as text, so force a Text number format for the assignment, then restore the
# cell style (no explicit style) to match the original formatting.
$numericLooking = @("D4", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D15", "D17", "D19", "D20", "D21", "D22", "D23", "D26", "D27", "D28", "D31", "D32", "D33", "D34", "D36", "D38", "D41", "D42", "D43", "D44", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $numericLooking) { $ws.Range($ref).NumberFormat = "@" }

$ws.Range("D4").Value = "1.002"
$ws.Range("D5").Value = "215.63"
$ws.Range("D6").Value = "0.5056"
$ws.Range("D8").Value = "0.2582"
$ws.Range("D9").Value = "0.06420"
$ws.Range("D10").Value = "20.35"
$ws.Range("D11").Value = "0.07791"
$ws.Range("D12").Value = "4.263"
$ws.Range("D15").Value = "0.5615"
$ws.Range("D17").Value = "63.25"
$ws.Range("D19").Value = "1.002"
$ws.Range("D20").Value = "192.88"
$ws.Range("D21").Value = "4.373"
$ws.Range("D22").Value = "9.912"
$ws.Range("D23").Value = "6.125"
$ws.Range("D26").Value = "140.96"
$ws.Range("D27").Value = "0.1236"
$ws.Range("D28").Value = "6.806"
$ws.Range("D31").Value = "0.04944"
$ws.Range("D32").Value = "3.284"
$ws.Range("D33").Value = "3.230"
$ws.Range("D34").Value = "1.568"
$ws.Range("D36").Value = "0.9037"
$ws.Range("D38").Value = "2.563"
$ws.Range("D41").Value = "0.9959"
$ws.Range("D42").Value = "5.479"
$ws.Range("D43").Value = "0.8020"
$ws.Range("D44").Value = "98.89"
$ws.Range("D47").Value = "55.60"
$ws.Range("D48").Value = "0.4273"
$ws.Range("D49").Value = "7.796"
$ws.Range("D50").Value = "0.05029"
$ws.Range("D51").Value = "1.0000"

foreach ($ref in $numericLooking) { $ws.Range($ref).Style = "Normal" }
